$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: OptionBurningMultiplier -> OptionFireMultiplier
$ws.Range("B11").Value = 40
$ws.Range("C11").Value = "OptionFireMultiplier"
$ws.Range("E11").Value = "1.2f"
$ws.Range("F11").Value = "DOT damage multiplier for fire attacks. While active, creature has Burning visual effect. 0.0x disables DOT from fire entirely."

# Row 12: OptionElectrocuteMultiplier -> OptionLightningMultiplier
$ws.Range("B12").Value = 50
$ws.Range("C12").Value = "OptionLightningMultiplier"
$ws.Range("E12").Value = "1.2f"
$ws.Range("F12").Value = "DOT damage multiplier for lightning attacks. While active, creature has Electrocute visual effect. 0.0x disables DOT from lightning entirely."
